# Fix the Japanese titles in the Double Moon Legend checklist so that the
# "japanese" column correctly matches the rulebook / supplement rows, and
# update the "english" column for the supplement row to match as well.
#
# Row 2 (1991, rulebook)   : japanese should read "...TRPGシステムブック"
# Row 3 (1992, supplement) : japanese should read "...スプリメントブック"
#                             english should read "Double Moon Legend Supplement Book"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for shared-string table layout: write the supplement
# strings first (english, then japanese), then the rulebook japanese string.
$ws.Range("C3").Value = "Double Moon Legend Supplement Book"
$ws.Range("B3").Value = "ダブルムーン伝説スプリメントブック"
$ws.Range("B2").Value = "ダブルムーン伝説TRPGシステムブック"

# Move the active selection to B4, matching the saved workbook state.
$ws.Range("B4").Select() | Out-Null
